$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '67.054.83'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value2 = '  -1.03%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '2.480.07'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value2 = '  -1.72%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value2 = '  -0.01%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '586.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value2 = '  -0.94%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '169.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value2 = '  -4.29%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value2 = '  +0.04%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.516'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value2 = '  -2.92%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '2.484.03'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value2 = '  -1.50%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value2 = '  -3.89%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value2 = '  +0.74%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.336'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value2 = '  -2.88%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '4.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value2 = '  -4.21%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '25.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value2 = '  -3.60%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '2.962.17'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value2 = '  -0.65%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value2 = '  -2.89%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '67.061.68'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value2 = '  -0.77%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '2.492.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value2 = '  -0.89%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '11.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value2 = '  +1.09%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '7.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value2 = '  -3.03%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '363.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value2 = '  +0.05%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '4.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value2 = '  -3.38%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '4.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value2 = '  -4.75%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value2 = 'Dai'
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '1.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value2 = '  +0.16%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value2 = 'Litecoin'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '71.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value2 = '  +0.27%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value2 = '  -6.45%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '9.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value2 = '  -8.08%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '1.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value2 = '  +8.45%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value2 = '  -0.81%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '0.0₃0925'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value2 = '  -6.37%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '8.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value2 = '  -1.94%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '514.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value2 = '  -4.91%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '1.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value2 = '  -1.81%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value2 = '  -5.38%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value2 = '  +0.07%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.126'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value2 = '  -2.43%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '156.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value2 = '  -0.01%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '19.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value2 = '  +0.96%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value2 = '  -2.95%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '18.56'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value2 = '  -0.53%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '1.74'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value2 = '  -2.91%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '4.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value2 = '  -4.61%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.331'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value2 = '  -7.26%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value2 = 'dogwifhat'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '2.46'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value2 = '  -1.92%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value2 = 'OKB'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '39.25'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value2 = '  -1.72%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value2 = 'Aave'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '142.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value2 = '  -2.44%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value2 = 'Filecoin'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '3.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value2 = '  -3.08%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value2 = 'ARBITRUM'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '0.536'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value2 = '  -4.46%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value2 = 'BabyDogeCoin'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '0.0₆0266'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value2 = '  -3.47%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value2 = 'Optimism'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '1.64'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value2 = '  -2.75%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value2 = 'Mantle'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '0.593'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value2 = '  -0.48%  '
$ws.Range('E51').Style = 'Normal'
